$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply centered alignment (same style already used on the trace table) to
#     the full header/body block before filling in the new columns, so the
#     new cells inherit it just like B5:D6 already have.
$ws.Range("B5:R5").HorizontalAlignment = -4108
$ws.Range("B5:R5").VerticalAlignment = -4108

# --- New trace-table columns (quicksort trace): i, j, pivot, loop conditions, etc.
$ws.Range("E5").Value = "i"
$ws.Range("F5").Value = "j"
$ws.Range("G5").Value = " pivot"
$ws.Range("H5").Value = "while(i<=j)"
$ws.Range("I5").Value = "while(numbers[i]<pivot)"
$ws.Range("J5").Value = "i++"
$ws.Range("K5").Value = "while(numbers[j] >pivot)"
$ws.Range("L5").Value = "j--"
$ws.Range("M5").Value = "i <=j"
$ws.Range("N5").Value = "exchange(i,j)"
$ws.Range("O5").Value = "low < j"
$ws.Range("P5").Value = "quicksort(low, j)"
$ws.Range("R5").Value = "quicksort( I, high)"
$ws.Range("Q5").Value = "I < high"

# --- Column widths for the newly added columns (best match achievable widths)
$ws.Columns("I").ColumnWidth = 22.5
$ws.Columns("K").ColumnWidth = 22.333333333333332
$ws.Columns("N").ColumnWidth = 15.166666666666666
$ws.Columns("P").ColumnWidth = 17.166666666666668
$ws.Columns("Q").ColumnWidth = 17.5
$ws.Columns("R").ColumnWidth = 21.833333333333332

# --- Selection moved to F10, matching the saved view state
$ws.Range("F10").Select() | Out-Null
